# Update crypto price/volume data per commit "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.353.51"
$ws.Range("E2").Value = "  +3.42%  "
$ws.Range("D3").Value = "3.059.55"
$ws.Range("E3").Value = "  +1.40%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "'549.47"
$ws.Range("E5").Value = "  +3.34%  "
$ws.Range("D6").Value = "'139.84"
$ws.Range("E6").Value = "  +3.93%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").Value = "3.055.03"
$ws.Range("E8").Value = "  +1.43%  "
$ws.Range("D9").Value = "'0.501"
$ws.Range("E9").Value = "  +0.55%  "
$ws.Range("D10").Value = "'6.39"
$ws.Range("E10").Value = "  +5.03%  "
$ws.Range("D11").Value = "'0.150"
$ws.Range("E11").Value = "  +0.70%  "
$ws.Range("D12").Value = "'0.452"
$ws.Range("E12").Value = "  +0.98%  "
$ws.Range("E13").Value = "  +2.62%  "
$ws.Range("D14").Value = "'34.73"
$ws.Range("E14").Value = "  +1.12%  "
$ws.Range("D15").Value = "3.549.04"
$ws.Range("E15").Value = "  +1.07%  "
$ws.Range("D16").Value = "63.313.93"
$ws.Range("E16").Value = "  +3.18%  "
$ws.Range("D17").Value = "3.052.45"
$ws.Range("E17").Value = "  +0.89%  "
$ws.Range("E18").Value = "  -1.46%  "
$ws.Range("D19").Value = "'6.73"
$ws.Range("E19").Value = "  +1.37%  "
$ws.Range("D20").Value = "'481.19"
$ws.Range("E20").Value = "  +3.00%  "
$ws.Range("D21").Value = "'13.65"
$ws.Range("E21").Value = "  +2.81%  "
$ws.Range("D22").Value = "'0.672"
$ws.Range("E22").Value = "  -1.25%  "
$ws.Range("D23").Value = "'7.19"
$ws.Range("E23").Value = "  +2.86%  "
$ws.Range("D24").Value = "'80.56"
$ws.Range("E24").Value = "  +1.28%  "
$ws.Range("D25").Value = "'12.48"
$ws.Range("E25").Value = "  +2.66%  "
$ws.Range("E26").Value = "  +0.16%  "
$ws.Range("D27").Value = "'2.75"
$ws.Range("E27").Value = "  +2.46%  "
$ws.Range("D28").Value = "'7.91"
$ws.Range("E28").Value = "  +0.66%  "
$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D29").Value = "'1.98"
$ws.Range("E29").Value = "  +3.99%  "
$ws.Range("B30").Value = "FirstDigitalUSD"
$ws.Range("C30").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D30").Value = "'0.997"
$ws.Range("E30").Value = "  -0.49%  "
$ws.Range("D31").Value = "'25.97"
$ws.Range("E31").Value = "  +1.05%  "
$ws.Range("D32").Value = "'1.15"
$ws.Range("E32").Value = "  +0.63%  "
$ws.Range("D33").Value = "'2.42"
$ws.Range("E33").Value = "  +5.70%  "
$ws.Range("D34").Value = "'5.68"
$ws.Range("E34").Value = "  +2.52%  "
$ws.Range("D35").Value = "'55.56"
$ws.Range("E35").Value = "  -0.38%  "
$ws.Range("D36").Value = "'5.96"
$ws.Range("E36").Value = "  +0.65%  "
$ws.Range("D37").Value = "'463.69"
$ws.Range("E37").Value = "  +0.31%  "
$ws.Range("D38").Value = "'0.0814"
$ws.Range("E38").Value = "  +3.04%  "
$ws.Range("D39").Value = "'0.0395"
$ws.Range("E39").Value = "  +2.21%  "
$ws.Range("D40").Value = "3.068.77"
$ws.Range("E40").Value = "  -5.01%  "
$ws.Range("E41").Value = "  -0.54%  "
$ws.Range("D42").Value = "'8.24"
$ws.Range("E42").Value = "  +0.75%  "
$ws.Range("D43").Value = "'2.59"
$ws.Range("E43").Value = "  +3.66%  "
$ws.Range("D44").Value = "'28.26"
$ws.Range("E44").Value = "  +2.04%  "
$ws.Range("D45").Value = "'0.253"
$ws.Range("E45").Value = "  +2.13%  "
$ws.Range("D47").Value = "'2.05"
$ws.Range("E47").Value = "  +2.00%  "
$ws.Range("E48").Value = "  +0.60%  "
$ws.Range("D49").Value = "'116.74"
$ws.Range("E49").Value = "  -2.13%  "
$ws.Range("D50").Value = "0.0₃0510"
$ws.Range("E50").Value = "  +2.66%  "
$ws.Range("E51").Value = "  +1.73%  "
